{"js": "// Update each arithmetic-problem answer cell in the table (100 replacements).\n// Each pair is [oldExpression, newExpression]; every old value is unique in\n// the document, so a literal (non-wildcard) search-and-replace is safe.\nconst pairs = [\n  [\"15+47=62\", \"91-56=35\"],\n  [\"15+56=71\", \"58+6=64\"],\n  [\"76-19=57\", \"35+6=41\"],\n  [\"52-3=49\", \"30-8=22\"],\n  [\"39+12=51\", \"58+33=91\"],\n  [\"6+15=21\", \"44-7=37\"],\n  [\"88-29=59\", \"15+77=92\"],\n  [\"56+6=62\", \"81-76=5\"],\n  [\"91-14=77\", \"5+89=94\"],\n  [\"83-36=47\", \"13+28=41\"],\n  [\"70-25=45\", \"63-36=27\"],\n  [\"63-14=49\", \"65-38=27\"],\n  [\"90-77=13\", \"23+69=92\"],\n  [\"13+38=51\", \"50-1=49\"],\n  [\"81-38=43\", \"79+5=84\"],\n  [\"74-46=28\", \"46-19=27\"],\n  [\"19+58=77\", \"73-27=46\"],\n  [\"67+25=92\", \"53-19=34\"],\n  [\"61-19=42\", \"66+6=72\"],\n  [\"8+8=16\", \"91-73=18\"],\n  [\"8+24=32\", \"84-48=36\"],\n  [\"12+59=71\", \"49+29=78\"],\n  [\"72-8=64\", \"44-15=29\"],\n  [\"45+28=73\", \"4+58=62\"],\n  [\"26-19=7\", \"94-66=28\"],\n  [\"81-69=12\", \"6+35=41\"],\n  [\"93-18=75\", \"52-25=27\"],\n  [\"18+54=72\", \"96-89=7\"],\n  [\"80-63=17\", \"85-49=36\"],\n  [\"18+38=56\", \"29+22=51\"],\n  [\"38+43=81\", \"63-15=48\"],\n  [\"53-18=35\", \"55+6=61\"],\n  [\"80-59=21\", \"50-31=19\"],\n  [\"22+39=61\", \"47-8=39\"],\n  [\"67+18=85\", \"25+18=43\"],\n  [\"9+77=86\", \"58+35=93\"],\n  [\"30-19=11\", \"80-76=4\"],\n  [\"92-4=88\", \"47+16=63\"],\n  [\"64-39=25\", \"2+39=41\"],\n  [\"6+85=91\", \"53-26=27\"],\n  [\"25+9=34\", \"54-16=38\"],\n  [\"48+43=91\", \"38+23=61\"],\n  [\"20-11=9\", \"9+54=63\"],\n  [\"9+5=14\", \"26+57=83\"],\n  [\"94-6=88\", \"83-38=45\"],\n  [\"77-38=39\", \"72-36=36\"],\n  [\"5+38=43\", \"22+69=91\"],\n  [\"69+24=93\", \"54+8=62\"],\n  [\"68+19=87\", \"18+43=61\"],\n  [\"27+57=84\", \"75+9=84\"],\n  [\"68+14=82\", \"59+7=66\"],\n  [\"92-37=55\", \"48+9=57\"],\n  [\"59+37=96\", \"58+13=71\"],\n  [\"35-9=26\", \"71-48=23\"],\n  [\"6+77=83\", \"92-58=34\"],\n  [\"47+15=62\", \"13-9=4\"],\n  [\"81-3=78\", \"59+35=94\"],\n  [\"65+9=74\", \"23+69=92\"],\n  [\"94-26=68\", \"93-29=64\"],\n  [\"72-3=69\", \"19+25=44\"],\n  [\"84-68=16\", \"29+2=31\"],\n  [\"81-29=52\", \"32-19=13\"],\n  [\"19+7=26\", \"42-29=13\"],\n  [\"94-58=36\", \"67+24=91\"],\n  [\"83-79=4\", \"13+49=62\"],\n  [\"25+57=82\", \"49+23=72\"],\n  [\"82-7=75\", \"76-48=28\"],\n  [\"44+18=62\", \"72-9=63\"],\n  [\"15+17=32\", \"76-49=27\"],\n  [\"74-45=29\", \"95-58=37\"],\n  [\"68+18=86\", \"30-25=5\"],\n  [\"79+2=81\", \"15+8=23\"],\n  [\"75-67=8\", \"52+39=91\"],\n  [\"24-17=7\", \"8+69=77\"],\n  [\"51-36=15\", \"69+7=76\"],\n  [\"81-15=66\", \"11-3=8\"],\n  [\"71-69=2\", \"67+27=94\"],\n  [\"77+19=96\", \"16+75=91\"],\n  [\"92-79=13\", \"7+44=51\"],\n  [\"78+5=83\", \"18+76=94\"],\n  [\"92-57=35\", \"77-39=38\"],\n  [\"25+8=33\", \"81-53=28\"],\n  [\"92-38=54\", \"61-57=4\"],\n  [\"19+43=62\", \"64-15=49\"],\n  [\"60-29=31\", \"72-54=18\"],\n  [\"28+46=74\", \"8+19=27\"],\n  [\"48+48=96\", \"86-59=27\"],\n  [\"68-29=39\", \"34+38=72\"],\n  [\"52-29=23\", \"29+32=61\"],\n  [\"45+26=71\", \"18+43=61\"],\n  [\"73-68=5\", \"88-9=79\"],\n  [\"83-69=14\", \"42-18=24\"],\n  [\"33-26=7\", \"71-6=65\"],\n  [\"69+16=85\", \"48+14=62\"],\n  [\"84+7=91\", \"45-9=36\"],\n  [\"95-36=59\", \"61-9=52\"],\n  [\"92-54=38\", \"3+78=81\"],\n  [\"44-38=6\", \"38+4=42\"],\n  [\"31-3=28\", \"92-19=73\"],\n  [\"62-3=59\", \"9+17=26\"]\n];\n\nconst body = context.document.body;\n\nfor (const [oldText, newText] of pairs) {\n  const results = body.search(oldText, { matchCase: true, matchWholeWord: false });\n  results.load(\"items\");\n  await context.sync();\n\n  for (let i = 0; i < results.items.length; i++) {\n    results.items[i].insertText(newText, Word.InsertLocation.replace);\n  }\n  await context.sync();\n}\n", "ps1": "# Update each arithmetic-problem answer cell in the table (100 replacements).\n$d = $word.ActiveDocument\n\n$pairs = @(\n    @(\"15+47=62\", \"91-56=35\"),\n    @(\"15+56=71\", \"58+6=64\"),\n    @(\"76-19=57\", \"35+6=41\"),\n    @(\"52-3=49\", \"30-8=22\"),\n    @(\"39+12=51\", \"58+33=91\"),\n    @(\"6+15=21\", \"44-7=37\"),\n    @(\"88-29=59\", \"15+77=92\"),\n    @(\"56+6=62\", \"81-76=5\"),\n    @(\"91-14=77\", \"5+89=94\"),\n    @(\"83-36=47\", \"13+28=41\"),\n    @(\"70-25=45\", \"63-36=27\"),\n    @(\"63-14=49\", \"65-38=27\"),\n    @(\"90-77=13\", \"23+69=92\"),\n    @(\"13+38=51\", \"50-1=49\"),\n    @(\"81-38=43\", \"79+5=84\"),\n    @(\"74-46=28\", \"46-19=27\"),\n    @(\"19+58=77\", \"73-27=46\"),\n    @(\"67+25=92\", \"53-19=34\"),\n    @(\"61-19=42\", \"66+6=72\"),\n    @(\"8+8=16\", \"91-73=18\"),\n    @(\"8+24=32\", \"84-48=36\"),\n    @(\"12+59=71\", \"49+29=78\"),\n    @(\"72-8=64\", \"44-15=29\"),\n    @(\"45+28=73\", \"4+58=62\"),\n    @(\"26-19=7\", \"94-66=28\"),\n    @(\"81-69=12\", \"6+35=41\"),\n    @(\"93-18=75\", \"52-25=27\"),\n    @(\"18+54=72\", \"96-89=7\"),\n    @(\"80-63=17\", \"85-49=36\"),\n    @(\"18+38=56\", \"29+22=51\"),\n    @(\"38+43=81\", \"63-15=48\"),\n    @(\"53-18=35\", \"55+6=61\"),\n    @(\"80-59=21\", \"50-31=19\"),\n    @(\"22+39=61\", \"47-8=39\"),\n    @(\"67+18=85\", \"25+18=43\"),\n    @(\"9+77=86\", \"58+35=93\"),\n    @(\"30-19=11\", \"80-76=4\"),\n    @(\"92-4=88\", \"47+16=63\"),\n    @(\"64-39=25\", \"2+39=41\"),\n    @(\"6+85=91\", \"53-26=27\"),\n    @(\"25+9=34\", \"54-16=38\"),\n    @(\"48+43=91\", \"38+23=61\"),\n    @(\"20-11=9\", \"9+54=63\"),\n    @(\"9+5=14\", \"26+57=83\"),\n    @(\"94-6=88\", \"83-38=45\"),\n    @(\"77-38=39\", \"72-36=36\"),\n    @(\"5+38=43\", \"22+69=91\"),\n    @(\"69+24=93\", \"54+8=62\"),\n    @(\"68+19=87\", \"18+43=61\"),\n    @(\"27+57=84\", \"75+9=84\"),\n    @(\"68+14=82\", \"59+7=66\"),\n    @(\"92-37=55\", \"48+9=57\"),\n    @(\"59+37=96\", \"58+13=71\"),\n    @(\"35-9=26\", \"71-48=23\"),\n    @(\"6+77=83\", \"92-58=34\"),\n    @(\"47+15=62\", \"13-9=4\"),\n    @(\"81-3=78\", \"59+35=94\"),\n    @(\"65+9=74\", \"23+69=92\"),\n    @(\"94-26=68\", \"93-29=64\"),\n    @(\"72-3=69\", \"19+25=44\"),\n    @(\"84-68=16\", \"29+2=31\"),\n    @(\"81-29=52\", \"32-19=13\"),\n    @(\"19+7=26\", \"42-29=13\"),\n    @(\"94-58=36\", \"67+24=91\"),\n    @(\"83-79=4\", \"13+49=62\"),\n    @(\"25+57=82\", \"49+23=72\"),\n    @(\"82-7=75\", \"76-48=28\"),\n    @(\"44+18=62\", \"72-9=63\"),\n    @(\"15+17=32\", \"76-49=27\"),\n    @(\"74-45=29\", \"95-58=37\"),\n    @(\"68+18=86\", \"30-25=5\"),\n    @(\"79+2=81\", \"15+8=23\"),\n    @(\"75-67=8\", \"52+39=91\"),\n    @(\"24-17=7\", \"8+69=77\"),\n    @(\"51-36=15\", \"69+7=76\"),\n    @(\"81-15=66\", \"11-3=8\"),\n    @(\"71-69=2\", \"67+27=94\"),\n    @(\"77+19=96\", \"16+75=91\"),\n    @(\"92-79=13\", \"7+44=51\"),\n    @(\"78+5=83\", \"18+76=94\"),\n    @(\"92-57=35\", \"77-39=38\"),\n    @(\"25+8=33\", \"81-53=28\"),\n    @(\"92-38=54\", \"61-57=4\"),\n    @(\"19+43=62\", \"64-15=49\"),\n    @(\"60-29=31\", \"72-54=18\"),\n    @(\"28+46=74\", \"8+19=27\"),\n    @(\"48+48=96\", \"86-59=27\"),\n    @(\"68-29=39\", \"34+38=72\"),\n    @(\"52-29=23\", \"29+32=61\"),\n    @(\"45+26=71\", \"18+43=61\"),\n    @(\"73-68=5\", \"88-9=79\"),\n    @(\"83-69=14\", \"42-18=24\"),\n    @(\"33-26=7\", \"71-6=65\"),\n    @(\"69+16=85\", \"48+14=62\"),\n    @(\"84+7=91\", \"45-9=36\"),\n    @(\"95-36=59\", \"61-9=52\"),\n    @(\"92-54=38\", \"3+78=81\"),\n    @(\"44-38=6\", \"38+4=42\"),\n    @(\"31-3=28\", \"92-19=73\"),\n    @(\"62-3=59\", \"9+17=26\"),\n)\n\nforeach ($pair in $pairs) {\n    $old = $pair[0]\n    $new = $pair[1]\n    $find = $d.Content.Find\n    $find.ClearFormatting()\n    $find.Text = $old\n    $find.Replacement.ClearFormatting()\n    $find.Replacement.Text = $new\n    $find.Execute($old, $false, $true, $false, $false, $false, $true, 0, $false, $new, 2) | Out-Null\n}\n"}
